# Update test data for Recommended Content:
# The cardImageSrc values (column J) on the "pages_with_recommended_content"
# sheet pointed at a dated image-style folder
# (".../cgov_image/featured/2019-11/"). Update them to the current,
# date-less folder (".../cgov_image/featured/").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pages_with_recommended_content")

$newCardImageSrc = "/sites/default/files/styles/cgov_featured/public/cgov_image/featured/"

$ws.Range("J2").Value = $newCardImageSrc
$ws.Range("J3").Value = $newCardImageSrc
$ws.Range("J4").Value = $newCardImageSrc

# Reflect the cursor/selection position left behind in the saved file.
$ws.Range("J7").Select()
